$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 8 (shifts rows 8-37 down to 9-38)
$ws.Rows.Item(8).Insert()

# Populate new row 8 with data
$ws.Cells.Item(8, 1).Value = 1
$ws.Cells.Item(8, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(8, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(8, 4).Value = 45030
$ws.Cells.Item(8, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(8, 5).Value = 15
$ws.Cells.Item(8, 6).Value = "Fruta"
$ws.Cells.Item(8, 7).Value = 100103
$ws.Cells.Item(8, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(8, 9).Value = 100103002
$ws.Cells.Item(8, 10).Value = "Ciruela"
$ws.Cells.Item(8, 11).Value = "Fortuna"
$ws.Cells.Item(8, 12).Value = "Segunda"
$ws.Cells.Item(8, 13).Value = 200
$ws.Cells.Item(8, 14).Value = 19000
$ws.Cells.Item(8, 15).Value = 20000
$ws.Cells.Item(8, 16).Value = 19500
$ws.Cells.Item(8, 17).Value = "`$/caja 18 kilos granel"
$ws.Cells.Item(8, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(8, 19).Value = 1083
$ws.Cells.Item(8, 20).Value = 18
